# Fix - Title of Entity in Menu not Work
# Rename sheets so the menu/tab titles show a user-friendly "Entity-Title"
# label, and restore the view state (active sheet / selected cells) that
# was captured when the workbook was last saved.

$wb = $excel.ActiveWorkbook

$wsRemarques       = $wb.Worksheets.Item(1)
$wsUser            = $wb.Worksheets.Item(2)
$wsArticle         = $wb.Worksheets.Item(3)
$wsOrder           = $wb.Worksheets.Item(4)
$wsDelivery        = $wb.Worksheets.Item(5)
$wsCategorie       = $wb.Worksheets.Item(6)
$wsProvider        = $wb.Worksheets.Item(7)
$wsOrderLine       = $wb.Worksheets.Item(8)

# Rename entity sheets to include their localized/display title.
$wsUser.Name      = "User-Utilisateur"
$wsDelivery.Name  = "Delivery-Laivraison"
$wsCategorie.Name = "CategorieArticle-CatégorieArtic"
$wsProvider.Name  = "Provider-Fournisseur"
$wsOrderLine.Name = "OrderLine-LigneCommande"

# Update the selections that moved on a couple of sheets.
$wsUser.Activate()
$wsUser.Range("D14").Select()

$wsArticle.Activate()
$wsArticle.Range("B18").Select()

# Delivery becomes the active/selected sheet (instead of OrderLine).
$wsDelivery.Activate()
$wsDelivery.Range("I9").Select()

# Scroll the sheet tab strip so "Order" is the first visible tab.
$excel.ActiveWindow.ScrollWorkbookTabs(3)
